# Add a new training-record row (row 69) to Sheet1's log table. The table
# rows share a common look (row height, per-column styles) so the new row
# is produced by duplicating the last row (68) - which carries that
# formatting, including the slightly different "K" column style used by
# the more recent entries - and then overwriting its cell values with the
# new training run's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 69

$ws.Rows.Item(68).Copy()
$ws.Rows.Item($row).Insert(-4121)  # xlShiftDown, carries row 68's formatting

# Values are written in the same order the original author entered them
# (not strictly left-to-right) so newly-minted shared-string entries land
# at the same indices as the canonical edit: D, E, K, L, then J. B/C keep
# the values copied from row 68 ("分类" / "14分类"), which already match.
$ws.Cells.Item($row, 1).Value = 43213.034722222219
$ws.Cells.Item($row, 4).Value = "batch_size=100 low_nums=5 use_biases=yes use_bn_low=True dropout_low=0.8"
$ws.Cells.Item($row, 5).Value = "最高标签，重新训练，PCA及Wavelet处理数据train-hjxh365-2018-4-16-day"
$ws.Cells.Item($row, 6).Value = 0.64
$ws.Cells.Item($row, 7).Value = 0.65
$ws.Cells.Item($row, 8).Value = 0.95
$ws.Cells.Item($row, 9).Value = 0.94
$ws.Cells.Item($row, 11).Value = "python feed_run.py --output_mode=classes --output_nodes=14 --input_nums=39 --input_nodes=39 --low_nums=5 --low_nodes=39 --low_fun=elu --use_bn_input=True --one_hot=True --input_fun=tanh --batch_size=100 --learning_rate=0.001 --train_mode=Adadelta --eval_size=5400 --test_size=1339 --use_biases=yes --dropout_low=0.8"
$ws.Cells.Item($row, 12).Value = "logs-hjxh-2018-4-23-class14-pca99-wavelet20-percent65"
$ws.Cells.Item($row, 10).Value = "经过约82小时，拟合精度还能提高，泛化精度看来还可以再提高一点，不过慢"

$ws.Rows.Item($row).RowHeight = 82.5

$ws.Range("J69").Select()
